# Update ValueSet-fr-editorial-status.xlsx metadata sheet
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# URL: https://hl7.fr/fhir/fr/medication/ValueSet/fr-editorial-status
#   -> https://hl7.fr/ig/fhir/medication/ValueSet/fr-editorial-status
$wsMeta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-editorial-status"

# Title: fix "InterOp'Sant" -> "Interop'Sant"
$wsMeta.Range("B5").Value = "value set Interop'Santé - Statut éditorial d'une valeur"

# Date: bump to new timestamp
$wsMeta.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction value: was empty, now FRANCE
$wsMeta.Range("B11").Value = "FRANCE"

$wsInc = $wb.Worksheets.Item("Include #0")

# System URI: https://hl7.fr/fhir/fr/medication/CodeSystem/fr-editorial-status
#   -> https://hl7.fr/ig/fhir/medication/CodeSystem/fr-editorial-status
$wsInc.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-editorial-status"
